$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row.
# A leading apostrophe forces Excel to treat a numeric-looking string (e.g. "1.00")
# as literal text so exact digits/trailing zeros are preserved instead of being
# parsed into a number.
$ws.Range("D2").Value = '60.276.39'
$ws.Range("E2").Value = '  -1.86%  '
$ws.Range("D3").Value = '3.389.46'
$ws.Range("E3").Value = '  -1.48%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''571.39'
$ws.Range("E5").Value = '  -1.39%  '
$ws.Range("D6").Value = '''141.23'
$ws.Range("E6").Value = '  -4.92%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '3.390.47'
$ws.Range("E8").Value = '  -1.48%  '
$ws.Range("D9").Value = '''0.474'
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E10").Value = '  -3.92%  '
$ws.Range("E11").Value = '  -0.86%  '
$ws.Range("D12").Value = '''0.393'
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").Value = '3.967.91'
$ws.Range("E13").Value = '  -1.50%  '
$ws.Range("D14").Value = '''28.19'
$ws.Range("E14").Value = '  +0.78%  '
$ws.Range("E15").Value = '  +1.00%  '
$ws.Range("E16").Value = '  -2.48%  '
$ws.Range("D17").Value = '3.391.37'
$ws.Range("E17").Value = '  -1.54%  '
$ws.Range("D18").Value = '60.454.54'
$ws.Range("E18").Value = '  -1.75%  '
$ws.Range("D19").Value = '''6.28'
$ws.Range("D20").Value = '''14.07'
$ws.Range("E20").Value = '  -2.06%  '
$ws.Range("D21").Value = '''9.12'
$ws.Range("E21").Value = '  -3.42%  '
$ws.Range("D22").Value = '''388.95'
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = '''0.562'
$ws.Range("E23").Value = '  -1.40%  '
$ws.Range("D24").Value = '''73.55'
$ws.Range("E24").Value = '  +1.21%  '
$ws.Range("E25").Value = '  +0.19%  '
$ws.Range("D26").Value = '''0.0000117'
$ws.Range("E26").Value = '  -4.17%  '
$ws.Range("D27").Value = '3.531.11'
$ws.Range("E27").Value = '  -1.67%  '
$ws.Range("D28").Value = '''0.180'
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").Value = '''7.39'
$ws.Range("E30").Value = '  -5.29%  '
$ws.Range("D31").Value = '''8.05'
$ws.Range("E31").Value = '  -2.36%  '
$ws.Range("E32").Value = '  -1.32%  '
$ws.Range("D33").Value = '''1.42'
$ws.Range("E33").Value = '  -7.54%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").Value = '''23.75'
$ws.Range("E35").Value = '  -0.93%  '
$ws.Range("E36").Value = '  -1.65%  '
$ws.Range("D37").Value = '3.419.54'
$ws.Range("E37").Value = '  -1.31%  '
$ws.Range("D38").Value = '''167.70'
$ws.Range("E38").Value = '  +0.97%  '
$ws.Range("E39").Value = '  -6.29%  '
$ws.Range("E40").Value = '  -4.23%  '
$ws.Range("D41").Value = '''0.0776'
$ws.Range("E41").Value = '  -1.95%  '
$ws.Range("D42").Value = '''27.15'
$ws.Range("E42").Value = '  +3.45%  '
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("D44").Value = '''1.00'
$ws.Range("E44").Value = '  -0.04%  '
$ws.Range("D45").Value = '''4.46'
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("E46").Value = '  -1.96%  '
$ws.Range("D47").Value = '''41.28'
$ws.Range("E47").Value = '  -2.38%  '
$ws.Range("D48").Value = '2.529.31'
$ws.Range("E48").Value = '  -3.18%  '
$ws.Range("E49").Value = '  -3.88%  '
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("D51").Value = '''6.84'
$ws.Range("E51").Value = '  -2.74%  '
